# Sync attendance_reports: normalize the "Recorded By" (column G) name
# ordering on the "Session Analysis Results" sheet.
#
# For every data row, the comma-separated list of names/emails in column G
# is reversed in place (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com"), except rows that include
# "admin@admin.com", which are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column G (Recorded By).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }
    if ($val -notlike "*,*") { continue }        # single entry, nothing to reorder
    if ($val -like "*admin@admin.com*") { continue }  # leave admin rows as-is

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $reversedParts = $parts[($parts.Count - 1)..0]
    $newVal = [string]::Join(", ", $reversedParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
